$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap "Periodo Mora" and "Valor Mora" values between row 16 and row 17
$ws.Range("E16").Value = "1811"
$ws.Range("E17").Value = "1810"
$ws.Range("F16").Value = 34200
$ws.Range("F17").Value = 38000
